$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 9285
$ws.Range("E2").Value = 710
$ws.Range("F2").Value = 718
$ws.Range("G2").Value = 745
$ws.Range("H2").Value = 537
$ws.Range("I2").Value = 307
$ws.Range("J2").Value = 230
$ws.Range("K2").Value = 11910
$ws.Range("L2").Value = 5239
$ws.Range("M2").Value = 6671
$ws.Range("N2").Value = 4468
$ws.Range("O2").Value = 2203
$ws.Range("P2").Value = 189
$ws.Range("Q2").Value = 1457
$ws.Range("R2").Value = -996
$ws.Range("S2").Value = 26
$ws.Range("T2").Value = 775
$ws.Range("U2").Value = 682
$ws.Range("V2").Value = 2207
$ws.Range("W2").Value = 7.65
$ws.Range("X2").Value = 5.79
$ws.Range("Y2").Value = 7.05
$ws.Range("Z2").Value = 5.22
$ws.Range("AA2").Value = 78.54000000000001
$ws.Range("AB2").Value = 2451.83
$ws.Range("AC2").Value = 811
$ws.Range("AD2").Value = 18.81
$ws.Range("AE2").Value = 12600
$ws.Range("AF2").Value = 1.21
$ws.Range("AG2").Value = 110
$ws.Range("AH2").Value = 0.72
$ws.Range("AI2").Value = 12.7
$ws.Range("AJ2").Value = 37882300

# Row 3
$ws.Range("D3").Value = 10908
$ws.Range("E3").Value = 912
$ws.Range("F3").Value = 912
$ws.Range("G3").Value = 876
$ws.Range("H3").Value = 629
$ws.Range("I3").Value = 279
$ws.Range("J3").Value = 350
$ws.Range("K3").Value = 12425
$ws.Range("L3").Value = 5471
$ws.Range("M3").Value = 6954
$ws.Range("N3").Value = 4703
$ws.Range("O3").Value = 2251
$ws.Range("P3").Value = 189
$ws.Range("Q3").Value = 866
$ws.Range("R3").Value = -1920
$ws.Range("S3").Value = 856
$ws.Range("T3").Value = 678
$ws.Range("U3").Value = 188
$ws.Range("V3").Value = 2395
$ws.Range("W3").Value = 8.369999999999999
$ws.Range("X3").Value = 5.77
$ws.Range("Y3").Value = 6.09
$ws.Range("Z3").Value = 5.17
$ws.Range("AA3").Value = 78.68000000000001
$ws.Range("AB3").Value = 2570.27
$ws.Range("AC3").Value = 737
$ws.Range("AD3").Value = 29.52
$ws.Range("AE3").Value = 13262
$ws.Range("AF3").Value = 1.64
$ws.Range("AG3").Value = 120
$ws.Range("AH3").Value = 0.55
$ws.Range("AI3").Value = 15.24
$ws.Range("AJ3").Value = 37882300

# Row 4
$ws.Range("D4").Value = 13464
$ws.Range("E4").Value = 1269
$ws.Range("F4").Value = 1269
$ws.Range("G4").Value = 1203
$ws.Range("H4").Value = 908
$ws.Range("I4").Value = 522
$ws.Range("J4").Value = 385
$ws.Range("K4").Value = 14083
$ws.Range("L4").Value = 6350
$ws.Range("M4").Value = 7732
$ws.Range("N4").Value = 5106
$ws.Range("O4").Value = 2627
$ws.Range("P4").Value = 189
$ws.Range("Q4").Value = 1396
$ws.Range("R4").Value = -253
$ws.Range("S4").Value = -3
$ws.Range("T4").Value = 656
$ws.Range("U4").Value = 740
$ws.Range("V4").Value = 2541
$ws.Range("W4").Value = 9.43
$ws.Range("X4").Value = 6.74
$ws.Range("Y4").Value = 10.65
$ws.Range("Z4").Value = 6.85
$ws.Range("AA4").Value = 82.13
$ws.Range("AB4").Value = 2812.27
$ws.Range("AC4").Value = 1379
$ws.Range("AD4").Value = 12.76
$ws.Range("AE4").Value = 14397
$ws.Range("AF4").Value = 1.22
$ws.Range("AG4").Value = 130
$ws.Range("AH4").Value = 0.74
$ws.Range("AI4").Value = 8.83
$ws.Range("AJ4").Value = 37882300

# Row 5
$ws.Range("D5").Value = 14501
$ws.Range("E5").Value = 1226
$ws.Range("F5").Value = 1226
$ws.Range("G5").Value = 1244
$ws.Range("H5").Value = 928
$ws.Range("I5").Value = 492
$ws.Range("J5").Value = 435
$ws.Range("K5").Value = 17424
$ws.Range("L5").Value = 8644
$ws.Range("M5").Value = 8779
$ws.Range("N5").Value = 5579
$ws.Range("O5").Value = 3201
$ws.Range("P5").Value = 189
$ws.Range("Q5").Value = 1289
$ws.Range("R5").Value = -242
$ws.Range("S5").Value = 414
$ws.Range("T5").Value = 785
$ws.Range("U5").Value = 504
$ws.Range("V5").Value = 3828
$ws.Range("W5").Value = 8.460000000000001
$ws.Range("X5").Value = 6.4
$ws.Range("Y5").Value = 9.210000000000001
$ws.Range("Z5").Value = 5.89
$ws.Range("AA5").Value = 98.45999999999999
$ws.Range("AB5").Value = 3025.16
$ws.Range("AC5").Value = 1299
$ws.Range("AD5").Value = 11.27
$ws.Range("AE5").Value = 15731
$ws.Range("AF5").Value = 0.93
$ws.Range("AG5").Value = 130
$ws.Range("AH5").Value = 0.89
$ws.Range("AI5").Value = 9.369999999999999
$ws.Range("AJ5").Value = 37882300

# Row 6
$ws.Range("D6").Value = 16119
$ws.Range("E6").Value = 1283
$ws.Range("F6").Value = 1283
$ws.Range("G6").Value = 1339
$ws.Range("H6").Value = 993
$ws.Range("I6").Value = 601
$ws.Range("K6").Value = 21730
$ws.Range("L6").Value = 11330
$ws.Range("M6").Value = 10400
$ws.Range("N6").Value = 6443
$ws.Range("P6").Value = 189
$ws.Range("Q6").Value = 2484
$ws.Range("R6").Value = -2644
$ws.Range("S6").Value = 272
$ws.Range("T6").Value = 1067
$ws.Range("U6").Value = 1417
$ws.Range("V6").Value = 5231
$ws.Range("W6").Value = 7.96
$ws.Range("X6").Value = 6.16
$ws.Range("Y6").Value = 10.01
$ws.Range("Z6").Value = 5.07
$ws.Range("AA6").Value = 108.94
$ws.Range("AB6").Value = 3316.55
$ws.Range("AC6").Value = 1588
$ws.Range("AD6").Value = 11.31
$ws.Range("AE6").Value = 17008
$ws.Range("AF6").Value = 1.06
$ws.Range("AG6").Value = 180
$ws.Range("AH6").Value = 1
$ws.Range("AI6").Value = 11.34
$ws.Range("AJ6").Value = 37882300

# Row 7
$ws.Range("D7").Value = 19735
$ws.Range("E7").Value = 1485
$ws.Range("G7").Value = 1613
$ws.Range("H7").Value = 1183
$ws.Range("I7").Value = 616
$ws.Range("K7").Value = 21855
$ws.Range("L7").Value = 10923
$ws.Range("M7").Value = 10932
$ws.Range("N7").Value = 6975
$ws.Range("P7").Value = 189
$ws.Range("Q7").Value = 529
$ws.Range("R7").Value = -1065
$ws.Range("S7").Value = -439
$ws.Range("T7").Value = 943
$ws.Range("W7").Value = 7.52
$ws.Range("X7").Value = 5.99
$ws.Range("Y7").Value = 9.18
$ws.Range("Z7").Value = 5.43
$ws.Range("AA7").Value = 99.92
$ws.Range("AC7").Value = 1626
$ws.Range("AD7").Value = 13.65
$ws.Range("AE7").Value = 18412
$ws.Range("AF7").Value = 1.21
$ws.Range("AG7").Value = 200
$ws.Range("AH7").Value = 0.9
$ws.Range("AI7").Value = 12.3

# Row 8
$ws.Range("D8").Value = 20931
$ws.Range("E8").Value = 1710
$ws.Range("G8").Value = 1837
$ws.Range("H8").Value = 1395
$ws.Range("I8").Value = 733
$ws.Range("K8").Value = 22182
$ws.Range("L8").Value = 10592
$ws.Range("M8").Value = 11590
$ws.Range("N8").Value = 7633
$ws.Range("P8").Value = 189
$ws.Range("Q8").Value = 1148
$ws.Range("R8").Value = -1025
$ws.Range("S8").Value = -409
$ws.Range("T8").Value = 1014
$ws.Range("W8").Value = 8.17
$ws.Range("X8").Value = 6.67
$ws.Range("Y8").Value = 10.04
$ws.Range("Z8").Value = 6.34
$ws.Range("AA8").Value = 91.39
$ws.Range("AC8").Value = 1935
$ws.Range("AD8").Value = 11.47
$ws.Range("AE8").Value = 20149
$ws.Range("AF8").Value = 1.1
$ws.Range("AG8").Value = 220
$ws.Range("AH8").Value = 0.99
$ws.Range("AI8").Value = 11.37

# Row 9
$ws.Range("D9").Value = 21886
$ws.Range("E9").Value = 1818
$ws.Range("G9").Value = 1944
$ws.Range("H9").Value = 1484
$ws.Range("I9").Value = 788
$ws.Range("K9").Value = 22582
$ws.Range("L9").Value = 10288
$ws.Range("M9").Value = 12294
$ws.Range("N9").Value = 8337
$ws.Range("P9").Value = 189
$ws.Range("Q9").Value = 1339
$ws.Range("R9").Value = -1116
$ws.Range("S9").Value = -383
$ws.Range("T9").Value = 1101
$ws.Range("W9").Value = 8.31
$ws.Range("X9").Value = 6.78
$ws.Range("Y9").Value = 9.869999999999999
$ws.Range("Z9").Value = 6.63
$ws.Range("AA9").Value = 83.68000000000001
$ws.Range("AC9").Value = 2080
$ws.Range("AD9").Value = 10.67
$ws.Range("AE9").Value = 22008
$ws.Range("AF9").Value = 1.01
$ws.Range("AG9").Value = 240
$ws.Range("AH9").Value = 1.08
$ws.Range("AI9").Value = 11.54

$ws.Range("U7").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("U9").ClearContents()